# Apply the "Clean up tests." edit:
#  - Insert a new ListParagraph bullet ("Clean up tests.") right before the
#    existing "Split tests even further..." bullet, at the same outline
#    level (ilvl=1, numId=1).
#  - Move the (hidden) "_GoBack" bookmark from the end of the last bullet
#    ("Re-add all components from v1.0.0.") to the end of the newly
#    inserted "Clean up tests." bullet, since that's where the edit
#    happened.

$d = $word.ActiveDocument

# --- Step 1: locate the "Split tests even further..." paragraph by its
#     text (robust against any pre-existing paragraph numbering) and
#     insert a brand new paragraph immediately before it, inheriting its
#     paragraph formatting (ListParagraph style, ilvl=1, numId=1).
$findRange = $d.Content
$findRange.Find.Execute("Split tests even further. E.g. Read -> ReadPod, ReadList, ReadString, etc.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitParaIndex = $findRange.Paragraphs(1).Index

$d.Paragraphs($splitParaIndex).Range.InsertParagraphBefore()

# --- Step 2: the freshly minted (empty) paragraph now sits right before
#     "Split tests even further...". Fill in its text.
$newPara = $d.Paragraphs($splitParaIndex)
$newRange = $newPara.Range.Duplicate
$newRange.MoveEnd(1, -1) | Out-Null
$newRange.Text = "Clean up tests."

# --- Step 3: remove the old hidden "_GoBack" bookmark (currently sitting
#     at the end of the very last bullet in the list).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 4: re-create "_GoBack" as a zero-length bookmark at the end of
#     the text we just typed ("Clean up tests."), matching where Word
#     leaves the insertion-point marker after an edit.
#     A temporary one-character pad works around a collapsed-range edge
#     case that otherwise misfires exactly at the end of a paragraph's
#     text; the pad is removed again once the bookmark is anchored.
$cleanPara = $d.Paragraphs($splitParaIndex)
$cleanTextRange = $cleanPara.Range.Duplicate
$cleanTextRange.MoveEnd(1, -1) | Out-Null
$endPos = $cleanTextRange.End

$pad = $d.Range($endPos, $endPos)
$pad.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padRange = $d.Range($endPos, $endPos + 1)
$padRange.Delete()
